$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 5: "MLE" (Male) was incorrectly labelled with the Arabic word for
# "Female" (أنثى). Correct it to the Arabic word for "Male" (الذكر). This
# adds a new shared string and repoints B5 at it.
$ws.Range("B5").Value() = "الذكر"

# Autofit column B (name column) to its contents.
$ws.Columns("B").AutoFit() | Out-Null

# Move/record the active selection as it was when the workbook was saved.
$ws.Range("D16").Select() | Out-Null

# Restore the page setup (paper size / orientation) recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
